$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New adjusted rent contribution value (row 42, column C = "041 Alquiler de vivienda")
$newC42 = 242.07

# Sum of column B for the detail rows (2..41), i.e. excluding the rent row (42)
$sumB = 0
for ($r = 2; $r -le 41; $r++) {
    $sumB += $ws.Cells.Item($r, 2).Value2
}

# Ratio used to rescale all other weightings so that, together with the new
# rent contribution, everything still sums to 1000
$ratio = (1000 - $newC42) / $sumB

for ($r = 2; $r -le 41; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $b * $ratio
}

$ws.Cells.Item(42, 3).Value = $newC42
